# TestIfStatement.docx: convert the old "@if / @endif" pseudo-syntax
# paragraphs into the new "{{if(...)if ... endif}}" template syntax,
# merging each if/body/endif trio of paragraphs into a single paragraph
# built out of many small runs (one run per literal token), and folding
# the "Points" if-block in between the two "Company_Name" if-blocks so
# the three surviving paragraphs hold:
#   1) MiniSofteware-if + Points-if + MaxiSoftware-if
#   2) CreateDate <= 2020-12-31 if
#   3) CreateDate > 2020-12-31 if

$d = $word.ActiveDocument

# Rebuilds paragraph number $paraIndex (1-based) so its content is made
# up of one run per entry of $pieces, in order. Works by writing the
# first piece directly into the paragraph, then for each remaining
# piece: splitting a fresh paragraph right after the current one,
# filling that new paragraph with the piece (its own run), and deleting
# the paragraph mark that separates them again -- merging two
# previously-distinct paragraphs always keeps their runs distinct,
# which is how we get the fine-grained run split the target markup
# wants instead of one big coalesced run.
function Set-ParagraphRuns($paraIndex, $pieces) {
    $p = $d.Paragraphs.Item($paraIndex)
    $s = $p.Range.Start
    $e = $p.Range.End - 1
    $r = $d.Range($s, $e)
    $r.Text = $pieces[0]

    for ($k = 1; $k -lt $pieces.Count; $k++) {
        $curP = $d.Paragraphs.Item($paraIndex)
        $endPos = $curP.Range.End - 1
        $ins = $d.Range($endPos, $endPos)
        $ins.InsertParagraphAfter()

        $newP = $d.Paragraphs.Item($paraIndex + 1)
        $ns = $newP.Range.Start
        $ne = $newP.Range.End - 1
        $nr = $d.Range($ns, $ne)
        $nr.Text = $pieces[$k]

        $curP2 = $d.Paragraphs.Item($paraIndex)
        $markStart = $curP2.Range.End - 1
        $markRange = $d.Range($markStart, $markStart + 1)
        $markRange.Delete()
    }
}

# Deletes $count whole paragraphs starting at $paraIndex (1-based),
# each removal (text + its paragraph mark) folding the next surviving
# paragraph up into this slot.
function Remove-Paragraphs($paraIndex, $count) {
    for ($i = 0; $i -lt $count; $i++) {
        $p = $d.Paragraphs.Item($paraIndex)
        $p.Range.Delete()
    }
}

$p1Runs = @("{{if", "(", "{{", "Company_Name", "}}", ",", "==", ",", "MiniSofteware", ")", "if", "First if chosen: {{Company_Name}}", "endif}}{{if", "(", "{{", "Points", "}}", ",>=,100)ifPoints are greater than 100", "endif}}{{if", "(", "{{", "Company_Name", "}}", ",", "==", ",", "MaxiSoftware", ")if", " ", "Second ", "if chosen: {{Company_Name}}", "endif}}")
$p2Runs = @("{{if", "(", "{{", "CreateDate", "}}", ",", "<=", ",", "2020-12-31", ")if", "CreateDate is not greater than 2021", "endif}}")
$p3Runs = @("{{if", "(", "{{", "CreateDate", "}}", ",>,", "2020-12-31", ")if", "CreateDate is not less than 2021", "endif}}")

# --- Paragraph 1: "@if Company_Name == MiniSofteware" becomes the
#     merged MiniSofteware-if / Points-if / MaxiSoftware-if paragraph.
Set-ParagraphRuns 1 $p1Runs

# The original paragraphs 2-9 ("First if chosen...", "@endif",
# "@if Company_Name == MaxiSoftware", "Second if chosen...", "@endif",
# "@if Points >= 100", "Points are greater than 100", "@endif") are now
# folded into paragraph 1's text above, so drop them.
Remove-Paragraphs 2 8

# --- Paragraph 2 (was "@if CreateDate <= 2020-12-31"): rebuild in place.
Set-ParagraphRuns 2 $p2Runs
Remove-Paragraphs 3 2

# --- Paragraph 3 (was "@if CreateDate > 2020-12-31"): rebuild in place.
Set-ParagraphRuns 3 $p3Runs
Remove-Paragraphs 4 2

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : $($d.Paragraphs.Item($i).Range.Text)"
}
